# forests-scraped.xlsx data refresh (2025-10-14 12:18)
#
# The 4 listings that were sitting in the "New" sheet (rows 2-5) have been
# triaged and move to the bottom of the "Previously added" sheet (becoming
# rows 176-179 there). Two freshly scraped listings take their place at the
# top of "New" (rows 2-3), so "New" shrinks from 5 rows to 3.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

# --- remember the 4 outgoing hyperlink targets (keyed by source row) before
#     anything on the "New" sheet is touched --------------------------------
$hlByRow = @{}
foreach ($hl in $ws2.Hyperlinks) {
    $hlByRow[$hl.Range.Row] = $hl.Address
}

# --- 1) move the 4 rows from "New" (A2:F5) to the end of "Previously added"
#        (A176:F179). Hyperlinks.Add is done *first* against the (currently
#        blank) destination cells; Range.Copy right after restores the
#        correct per-cell styles (s="3"/"4"/"2") that Hyperlinks.Add would
#        otherwise overwrite with its own default "Hyperlink" style, while
#        the hyperlink relationship itself survives the subsequent copy ----
for ($i = 0; $i -le 3; $i++) {
    $srcRow = 2 + $i
    $dstRow = 176 + $i
    $ws1.Hyperlinks.Add($ws1.Range("A" + $dstRow), $hlByRow[$srcRow]) | Out-Null
}

$ws2.Range("A2:F5").Copy($ws1.Range("A176:F179"))

# --- 2) clear the "New" sheet's old hyperlinks, then drop the now-migrated
#        rows 4:5 (rows 2:3 are left in place - with their original styling -
#        to be overwritten with the 2 newly scraped listings) --------------
$ws2.Hyperlinks.Delete() | Out-Null
$ws2.Rows("4:5").Delete() | Out-Null

# --- 3) write the 2 newly scraped listings into "New" A2:F3. As above,
#        Hyperlinks.Add first stamps its own style on column A, so a 1-cell
#        Range.Copy from a still correctly-styled donor cell (row 175 on
#        "Previously added") restores s="3" before the real URL text (and
#        the rest of the row's values) are written ---------------------
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/naujenes-pag/achid.html") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://www.ss.com/msg/lv/real-estate/wood/jekabpils-and-reg/vipes-pag/odlhl.html") | Out-Null

$ws1.Range("A175").Copy($ws2.Range("A2"))
$ws1.Range("A175").Copy($ws2.Range("A3"))

$ws2.Cells.Item(2, 1).Value = "https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/naujenes-pag/achid.html"
$ws2.Cells.Item(2, 2).Value = "4 500 €"
$ws2.Cells.Item(2, 3).Value = "Daugavpils un raj."
$ws2.Cells.Item(2, 4).Value = "1 ha."
$ws2.Cells.Item(2, 5).Value = ""
$ws2.Cells.Item(2, 6).Value = 45943.87986111111

$ws2.Cells.Item(3, 1).Value = "https://www.ss.com/msg/lv/real-estate/wood/jekabpils-and-reg/vipes-pag/odlhl.html"
$ws2.Cells.Item(3, 2).Value = "35 000 €"
$ws2.Cells.Item(3, 3).Value = "Jēkabpils un raj."
$ws2.Cells.Item(3, 4).Value = "5 ha."
$ws2.Cells.Item(3, 5).Value = "56960010027"
$ws2.Cells.Item(3, 6).Value = 45944.61388888889
